$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 5 new columns before column D (the existing quarterly
#    data in D:H - the most recent 5 quarters - shifts right to I:M,
#    making room for 5 older quarters in the new D:H).
# ------------------------------------------------------------------
$ws.Range("D1:H1").EntireColumn.Insert()

# ------------------------------------------------------------------
# 2. Re-apply explicit column widths so both the new (D:H) and the
#    shifted (I:M) quarter columns keep the original look
#    (29 / 31 / 29 / 29 / 29 chars - matching the pre-existing D:H
#    pattern used for every quarter column).
# ------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 28.17
$ws.Columns("E").ColumnWidth = 28.17
$ws.Columns("F").ColumnWidth = 30.17
$ws.Columns("G").ColumnWidth = 28.17
$ws.Columns("H").ColumnWidth = 28.17
$ws.Columns("I").ColumnWidth = 28.17
$ws.Columns("J").ColumnWidth = 30.17
$ws.Columns("K").ColumnWidth = 28.17
$ws.Columns("L").ColumnWidth = 28.17
$ws.Columns("M").ColumnWidth = 28.17

# --- Header row 8: new (earlier) quarter period labels ---
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"

# --- Header row 9: new publish-date labels ---
$ws.Range("D9").Value = "1400-08-25 (4)"
$ws.Range("E9").Value = "1400-10-29 (2)"
$ws.Range("F9").Value = "1401-04-01 (8)"
$ws.Range("G9").Value = "1401-04-28 (2)"
$ws.Range("H9").Value = "1401-08-25 (4)"

# --- Data rows 11-27: new (earlier) quarter figures ---
$ws.Range("D11").Value = 833625
$ws.Range("E11").Value = 898242
$ws.Range("F11").Value = 1227525
$ws.Range("G11").Value = 1106328
$ws.Range("H11").Value = 1406758
$ws.Range("D12").Value = -653450
$ws.Range("E12").Value = -717180
$ws.Range("F12").Value = -925967
$ws.Range("G12").Value = -840391
$ws.Range("H12").Value = -1064338
$ws.Range("D13").Value = 180175
$ws.Range("E13").Value = 181062
$ws.Range("F13").Value = 301558
$ws.Range("G13").Value = 265937
$ws.Range("H13").Value = 342420
$ws.Range("D14").Value = -39276
$ws.Range("E14").Value = -46520
$ws.Range("F14").Value = -58880
$ws.Range("G14").Value = -60644
$ws.Range("H14").Value = -57447
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("D16").Value = -7205
$ws.Range("E16").Value = -7335
$ws.Range("F16").Value = -4393
$ws.Range("G16").Value = -5225
$ws.Range("H16").Value = -4686
$ws.Range("D17").Value = 133694
$ws.Range("E17").Value = 127207
$ws.Range("F17").Value = 238285
$ws.Range("G17").Value = 200068
$ws.Range("H17").Value = 280287
$ws.Range("D18").Value = -306
$ws.Range("E18").Value = -1759
$ws.Range("F18").Value = -2425
$ws.Range("G18").Value = -2138
$ws.Range("H18").Value = -1825
$ws.Range("D19").Value = 5465
$ws.Range("E19").Value = 304
$ws.Range("F19").Value = 5866
$ws.Range("G19").Value = 576
$ws.Range("H19").Value = 861
$ws.Range("D20").Value = 138853
$ws.Range("E20").Value = 125752
$ws.Range("F20").Value = 241726
$ws.Range("G20").Value = 198506
$ws.Range("H20").Value = 279323
$ws.Range("D21").Value = -22733
$ws.Range("E21").Value = -23983
$ws.Range("F21").Value = -34808
$ws.Range("G21").Value = -37223
$ws.Range("H21").Value = -50151
$ws.Range("D22").Value = 116120
$ws.Range("E22").Value = 101769
$ws.Range("F22").Value = 206918
$ws.Range("G22").Value = 161283
$ws.Range("H22").Value = 229172
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("D24").Value = 116120
$ws.Range("E24").Value = 101769
$ws.Range("F24").Value = 206918
$ws.Range("G24").Value = 161283
$ws.Range("H24").Value = 229172
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = 23
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 51
$ws.Range("D26").Value = 4484000
$ws.Range("E26").Value = 4484000
$ws.Range("F26").Value = 4484000
$ws.Range("G26").Value = 4484000
$ws.Range("H26").Value = 4484000
$ws.Range("D27").Value = 26
$ws.Range("E27").Value = 23
$ws.Range("F27").Value = 46
$ws.Range("G27").Value = 36
$ws.Range("H27").Value = 51


# ------------------------------------------------------------------
# 3. Scroll the view so the newly-added (older) quarters are visible
#    next to the most recent ones.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 8
